$d = $word.ActiveDocument

# --- Merge split runs in title/author/abstract paragraphs into single runs ---
function Set-ParaXML($paraIndex, $styleVal, $text) {
    $p = $d.Paragraphs($paraIndex).Range
    $full = $d.Range($p.Start, $p.End - 1)
    $escText = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="' + $styleVal + '"/></w:pPr><w:r><w:t xml:space="preserve">' + $escText + '</w:t></w:r></w:p>'
    $full.InsertXML($xmlFrag)
}

Set-ParaXML 1 "Title" "Answers: Introduction to integration"
Set-ParaXML 2 "Author" "Donald Campbell"
Set-ParaXML 4 "Abstract" "Answers to questions relating to the guide on introduction to integration."

# --- Fix m:dPr child-element order: sepChr must come before endChr ---
$oMaths = $d.OMaths
$oMaths.Item(18).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>4</m:t></m:r></m:num><m:den><m:r><m:t>x</m:t></m:r></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="|" /><m:sepChr m:val="" /><m:endChr m:val="|" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(19).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="|" /><m:sepChr m:val="" /><m:endChr m:val="|" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(20).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(21).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(22).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>5</m:t></m:r></m:num><m:den><m:r><m:t>6</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(23).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(24).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∫</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>x</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>C</m:t></m:r></m:oMath>')
$oMaths.Item(29).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>0</m:t></m:r></m:sub><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:d></m:sup><m:e><m:r><m:t>4</m:t></m:r></m:e></m:nary><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:t>x</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>8</m:t></m:r></m:oMath>')
$oMaths.Item(30).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>0</m:t></m:r></m:sub><m:sup><m:r><m:t>5</m:t></m:r></m:sup><m:e><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:sup></m:sSup></m:e></m:nary><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>15</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:oMath>')
$oMaths.Item(31).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>1</m:t></m:r></m:sub><m:sup><m:r><m:t>2</m:t></m:r></m:sup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r></m:e></m:nary><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r></m:sup></m:sSup><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:oMath>')
$oMaths.Item(32).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>1</m:t></m:r></m:sub><m:sup><m:r><m:t>2</m:t></m:r></m:sup><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r></m:num><m:den><m:r><m:t>x</m:t></m:r></m:den></m:f></m:e></m:nary><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:d></m:oMath>')
$oMaths.Item(35).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>0</m:t></m:r></m:sub><m:sup><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>2</m:t></m:r></m:sup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e></m:nary><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>1</m:t></m:r></m:oMath>')
$oMaths.Item(36).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>0</m:t></m:r></m:sub><m:sup><m:r><m:t>π</m:t></m:r></m:sup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e></m:nary><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>0</m:t></m:r></m:oMath>')
$oMaths.Item(37).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>0</m:t></m:r></m:sub><m:sup><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>4</m:t></m:r></m:sup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e></m:nary><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:oMath>')
$oMaths.Item(38).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:t>0</m:t></m:r></m:sub><m:sup><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:sup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e></m:nary><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:rad><m:radPr><m:degHide m:val="on" /></m:radPr><m:deg /><m:e><m:r><m:t>3</m:t></m:r></m:e></m:rad></m:num><m:den><m:r><m:t>4</m:t></m:r></m:den></m:f></m:oMath>')
$oMaths.Item(39).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:nary><m:naryPr><m:chr m:val="∫" /><m:limLoc m:val="subSup" /><m:subHide m:val="off" /><m:supHide m:val="off" /></m:naryPr><m:sub><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>4</m:t></m:r></m:sub><m:sup><m:r><m:t>0</m:t></m:r></m:sup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e></m:nary><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>d</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r><m:rad><m:radPr><m:degHide m:val="on" /></m:radPr><m:deg /><m:e><m:r><m:t>2</m:t></m:r></m:e></m:rad></m:den></m:f></m:oMath>')

Write-Output "done"
